$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.148.22"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.444.76"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'582.92"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'142.99"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.530"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "2.439.67"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +2.80%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").Value = "'26.45"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "2.886.64"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "62.044.41"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "2.430.29"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "'10.77"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'326.62"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -5.31%  "
$ws.Range("D25").Value = "'65.69"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "'9.10"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("D27").Value = "'602.05"
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("D28").Value = "0.0₃0964"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'7.99"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.41"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").Value = "'0.136"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "'152.85"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "'43.13"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "'142.04"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "0.0₆0266"
$ws.Range("E48").Value = "  +19.02%  "
$ws.Range("D49").Value = "'0.600"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'19.78"
$ws.Range("E51").Value = "  -0.39%  "
